$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = "2025-08-17 16:16:00"
$ws.Range("B15").Value = "create-team"
$ws.Range("C15").Value = "new-organization97"
$ws.Range("D15").Value = "devteam"
$ws.Range("I15").Value = "'False"
$ws.Range("I15").Style = "Normal"
